$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 989692.4
$ws.Range("J17").Value = 1011656.44
$ws.Range("L17").Value = 3034969.32
$ws.Range("N17").Value = -3035305.32
$ws.Range("H40").Value = 2733.75
$ws.Range("J40").Value = 3150.75
$ws.Range("L40").Value = 3150.75
$ws.Range("N40").Value = -3500.75
$ws.Range("H94").Value = 31325954
$ws.Range("I94").Value = 38462692
$ws.Range("J94").Value = 400095.34
$ws.Range("K94").Value = 38462692
$ws.Range("L94").Value = 400095.34
$ws.Range("M94").Value = -38462241
$ws.Range("N94").Value = -400997.34
$ws.Range("H132").Value = 2594.322
$ws.Range("I132").Value = 2417.2036
$ws.Range("K132").Value = 7251.610799999999
$ws.Range("M132").Value = -4721.610799999999
$ws.Range("H137").Value = 8144.8374
$ws.Range("I137").Value = 12659.667
$ws.Range("K137").Value = 37979.001
$ws.Range("M137").Value = -35429.001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2366.6943
$ws.Range("I2").Value = 2342.48
$ws.Range("J2").Value = 2421.7273
$ws.Range("K2").Value = 2342.48
$ws.Range("L2").Value = 2421.7273
$ws.Range("M2").Value = -2229.48
$ws.Range("N2").Value = -2647.7273
$ws.Range("H32").Value = 6490.946
$ws.Range("I32").Value = 6201.831
$ws.Range("K32").Value = 6201.831
$ws.Range("M32").Value = -5914.831
$ws.Range("H45").Value = 130539.5
$ws.Range("I45").Value = 171034.83
$ws.Range("K45").Value = 171034.83
$ws.Range("M45").Value = -170657.83
$ws.Range("H116").Value = 2366.6943
$ws.Range("I116").Value = 2342.48
$ws.Range("J116").Value = 2421.7273
$ws.Range("K116").Value = 2342.48
$ws.Range("L116").Value = 2421.7273
$ws.Range("M116").Value = -48.48000000000002
$ws.Range("N116").Value = -7009.7273
$ws.Range("H122").Value = 1159880.6
$ws.Range("I122").Value = 5947
$ws.Range("K122").Value = 17841
$ws.Range("M122").Value = -15391

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2366.6943
$ws.Range("I3").Value = 2342.48
$ws.Range("J3").Value = 2421.7273
$ws.Range("K3").Value = 2342.48
$ws.Range("L3").Value = 2421.7273
$ws.Range("M3").Value = -2228.48
$ws.Range("N3").Value = -2649.7273
$ws.Range("H19").Value = 45000
$ws.Range("J19").Value = 45000
$ws.Range("L19").Value = 45000
$ws.Range("H86").Value = 7911.1333
$ws.Range("I86").Value = 8693.166999999999
$ws.Range("K86").Value = 8693.166999999999
$ws.Range("M86").Value = -7570.166999999999
$ws.Range("H89").Value = 7911.1333
$ws.Range("I89").Value = 8693.166999999999
$ws.Range("K89").Value = 43465.835
$ws.Range("M89").Value = -37849.835
$ws.Range("H137").Value = 55000
$ws.Range("J137").Value = 55000
$ws.Range("L137").Value = 55000
$ws.Range("N137").Value = -65200
$ws.Range("N19").Value = -45346

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6057.7236
$ws.Range("I31").Value = 6382.2188
$ws.Range("J31").Value = 5365.467
$ws.Range("K31").Value = 6382.2188
$ws.Range("L31").Value = 5365.467
$ws.Range("M31").Value = -6087.2188
$ws.Range("N31").Value = -5955.467
$ws.Range("H34").Value = 6057.7236
$ws.Range("I34").Value = 6382.2188
$ws.Range("J34").Value = 5365.467
$ws.Range("K34").Value = 6382.2188
$ws.Range("L34").Value = 5365.467
$ws.Range("M34").Value = -6180.2188
$ws.Range("N34").Value = -5769.467
$ws.Range("H95").Value = 124821250
$ws.Range("J95").Value = 124821250
$ws.Range("L95").Value = 124821250
$ws.Range("N95").Value = -124826742
$ws.Range("H122").Value = 9200.4
$ws.Range("I122").Value = 10988.083
$ws.Range("K122").Value = 32964.249
$ws.Range("M122").Value = -30514.249
$ws.Range("H134").Value = 4328.0586
$ws.Range("I134").Value = 4922.793
$ws.Range("J134").Value = 878.6
$ws.Range("K134").Value = 14768.379
$ws.Range("L134").Value = 2635.8
$ws.Range("M134").Value = -12233.379
$ws.Range("N134").Value = -7705.8

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 528102.7
$ws.Range("I5").Value = 632.4
$ws.Range("J5").Value = 716484.9399999999
$ws.Range("K5").Value = 1897.2
$ws.Range("L5").Value = 2149454.82
$ws.Range("M5").Value = -1785.2
$ws.Range("N5").Value = -2149678.82
$ws.Range("H37").Value = 63729.21
$ws.Range("J37").Value = 63729.21
$ws.Range("L37").Value = 191187.63
$ws.Range("N37").Value = -191411.63
$ws.Range("H40").Value = 204.57143
$ws.Range("I40").Value = 222.5
$ws.Range("K40").Value = 890
$ws.Range("M40").Value = -821
$ws.Range("H124").Value = 10741.556
$ws.Range("I124").Value = 2000
$ws.Range("J124").Value = 11834.25
$ws.Range("K124").Value = 6000
$ws.Range("L124").Value = 35502.75
$ws.Range("M124").Value = -1090
$ws.Range("N124").Value = -45322.75
$ws.Range("H133").Value = 7412.3335
$ws.Range("I133").Value = 7412.3335
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 22237.0005
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -17177.0005
$ws.Range("H135").Value = 528102.7
$ws.Range("I135").Value = 632.4
$ws.Range("J135").Value = 716484.9399999999
$ws.Range("K135").Value = 5691.599999999999
$ws.Range("L135").Value = 6448364.459999999
$ws.Range("M135").Value = -3156.599999999999
$ws.Range("N135").Value = -6453434.459999999
$ws.Range("N133").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 50000
$ws.Range("J39").Value = 50000
$ws.Range("L39").Value = 50000
$ws.Range("N39").Value = -51064
$ws.Range("H92").Value = 25019.834
$ws.Range("J92").Value = 27024
$ws.Range("L92").Value = 27024
$ws.Range("N92").Value = -30768
$ws.Range("H122").Value = 9772.852000000001
$ws.Range("I122").Value = 7539.294
$ws.Range("J122").Value = 13569.9
$ws.Range("K122").Value = 22617.882
$ws.Range("L122").Value = 40709.7
$ws.Range("M122").Value = -20167.882
$ws.Range("N122").Value = -45609.7
$ws.Range("H126").Value = 11454.863
$ws.Range("I126").Value = 12545.546
$ws.Range("K126").Value = 37636.638
$ws.Range("M126").Value = -35166.638

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14444.85
$ws.Range("J7").Value = 6882.5835
$ws.Range("L7").Value = 6882.5835
$ws.Range("N7").Value = -7106.5835
$ws.Range("H126").Value = 14444.85
$ws.Range("J126").Value = 6882.5835
$ws.Range("L126").Value = 20647.7505
$ws.Range("N126").Value = -25587.7505
